# "change predict income algorithm" - roll the quarterly window forward by
# one quarter: drop the oldest quarter (column E) and append a new quarter
# (column N) both for the period-label headers and for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- 1. Shift the quarter-label headers in row 8 and row 24 one column to
#        the left, and put the new quarter label in column N.
$newQuarterLabel = "فصل چهارم منتهی به 1401/12"

foreach ($headerRow in 8, 24) {
    $labels = @()
    for ($c = 6; $c -le 14; $c++) {
        $labels += $ws.Cells.Item($headerRow, $c).Value2
    }
    $labels += $newQuarterLabel

    $col = 5
    foreach ($lbl in $labels) {
        $ws.Cells.Item($headerRow, $col).Value2 = $lbl
        $col++
    }
}

# --- 2. Shift the data rows the same way: columns E:N (5..14) each move one
#        column to the left, and the newest figure lands in column N (14).
#        Rows 11, 13 and 18 are all zeroes on both sides of the shift, so
#        their "new" quarter value is simply 0.
$newData = @{
    10 = 61533
    11 = 0
    12 = 229733
    13 = 0
    14 = 1932
    15 = 145
    16 = 3470
    17 = 40421
    18 = 0
    19 = 22802
    20 = 360036
    26 = 258
    27 = 508
}

$dataRows = 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 26, 27

foreach ($r in $dataRows) {
    $vals = @()
    for ($c = 6; $c -le 14; $c++) {
        $vals += $ws.Cells.Item($r, $c).Value2
    }
    $vals += $newData[$r]

    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value2 = $v
        $col++
    }
}
